$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.635.67"
$ws.Range("E2").Value = "  -0.48%  "

$ws.Range("D3").Value = "1.616.65"
$ws.Range("E3").Value = "  -0.85%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.991"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.58%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "209.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.13%  "

$ws.Range("E6").Value = "  -1.14%  "

$ws.Range("E7").Value = "  -0.63%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.04"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.82%  "

$ws.Range("E9").Value = "  -1.37%  "

$ws.Range("E10").Value = "  -1.35%  "

$ws.Range("E11").Value = "  -0.80%  "

$ws.Range("D12").Value = "1.846.81"
$ws.Range("E12").Value = "  -0.85%  "

$ws.Range("D13").Value = "1.634.49"
$ws.Range("E13").Value = "  +0.20%  "

$ws.Range("E14").Value = "  -1.68%  "

$ws.Range("E15").Value = "  -1.37%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.64"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.01%  "

$ws.Range("D17").Value = "27.652.43"
$ws.Range("E17").Value = "  -0.56%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "227.25"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.33%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.61"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.26%  "

$ws.Range("E20").Value = "  -1.20%  "

$ws.Range("E21").Value = "  -0.73%  "

$ws.Range("E22").Value = "  -1.33%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.02"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.32%  "

$ws.Range("E24").Value = "  -1.32%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.48"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.01%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.06%  "

$ws.Range("E27").Value = "  -0.77%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.40"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.47%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.991"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.69%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.17"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.12%  "

$ws.Range("E31").Value = "  -0.77%  "

$ws.Range("E32").Value = "  -1.72%  "

$ws.Range("E33").Value = "  -0.54%  "

$ws.Range("D34").Value = "1.388.77"
$ws.Range("E34").Value = "  -1.27%  "

$ws.Range("E35").Value = "  +1.00%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.995"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.63%  "

$ws.Range("E37").Value = "  -1.49%  "

$ws.Range("E38").Value = "  +0.29%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.555"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.78%  "

$ws.Range("E40").Value = "  -3.20%  "

$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.990"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.75%  "

$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.54%  "

$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.82"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.64%  "

$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "65.41"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.78%  "

$ws.Range("E45").Value = "  -2.93%  "

$ws.Range("D46").Value = "1.755.72"
$ws.Range("E46").Value = "  -0.97%  "

$ws.Range("E47").Value = "  -3.79%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.64"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.33%  "

$ws.Range("E49").Value = "  +0.94%  "

$ws.Range("E50").Value = "  -0.88%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.55"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.92%  "
